$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17-20: the remark column moves from "Ongoing" to "Done" now that this
# iteration is ready to be sent to the client.
$ws.Range("E17").Value = "Done"
$ws.Range("E18").Value = "Done"
$ws.Range("E19").Value = "Done"
$ws.Range("E20").Value = "Done"

# Row 21: add the new schedule entry (Sr. # 15). Copy the formatting from the
# row above (row 20) for the columns that carry a distinct per-row style
# (Date, Remarks, Hours), then fill in the new content.
$ws.Range("D20").Copy($ws.Range("D21"))
$ws.Range("E20").Copy($ws.Range("E21"))
$ws.Range("F20").Copy($ws.Range("F21"))

$ws.Range("B21").Value = 15
$ws.Range("C21").Value = "Preparation for actual data input"
$ws.Range("D21").Value = "'08 - 12 - 2019"
$ws.Range("E21").Value = "Ongoing"
$ws.Range("F21").Value = 4

$excel.CutCopyMode = 0

# Move the view down so the newly-added row is visible, matching the sheet
# being scrolled/selected around the new entry before sending to the client.
$ws.Range("A21").Select()
$ws.Application.ActiveWindow.ScrollRow = 16
